# "handle empty multiform sheet and add resistance to cell values"
#
# The workbook gains a 4th worksheet ("2012") mirroring the existing
# "2007"/"2008"/"2010" multiform sheets, but left empty aside from a
# single informational cell. It becomes the active sheet (so the
# previously-active "2010" sheet is no longer the selected tab), and its
# text is added to the shared-string table.

$wb = $excel.ActiveWorkbook

# Add the new sheet after the last existing sheet so the tab order stays
# 2007, 2008, 2010, 2012.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "2012"

# Only content: a single cell noting there is no instruction for this
# (empty) multiform revision.
$newSheet.Range("A1").Value = "aucune instruction"

# Match the selection left behind on the new sheet and make it the
# active / selected tab (deselecting "2010").
$newSheet.Range("D11").Select() | Out-Null
$newSheet.Activate() | Out-Null
